$d = $word.ActiveDocument

# The document currently ends with:
#   ... "- Cung cấp chức năng tìm kiếm để người dùng có thể tìm kiếm
#        sản phẩm nhanh chóng và dễ dàng. "      (paragraph N-2)
#   ""                                            (blank paragraph N-1)
#   "Bản đặc tả này cung cấp các chức năng ..."   (paragraph N, last)
#
# We need to remove the blank paragraph and the trailing "Bản đặc tả ..."
# paragraph entirely, so the document ends right after "...dễ dàng. ",
# immediately followed by the section properties.

# Step 1: merge away the blank paragraph by deleting its own paragraph
# mark. This merges it into the following ("Bản đặc tả...") paragraph
# without disturbing the formatting/identity of the preceding text
# paragraph we want to keep.
$n = $d.Paragraphs.Count
$blank = $d.Paragraphs.Item($n - 1)
$d.Range($blank.Range.End - 1, $blank.Range.End).Delete()

# Step 2: the "Bản đặc tả..." paragraph is now the last paragraph in the
# document; delete it entirely (text + its own paragraph mark).
$n = $d.Paragraphs.Count
$d.Paragraphs.Item($n).Range.Delete()
